$wb = $excel.ActiveWorkbook

# --- 1. Insert new "2022-Q1" sheet right before the "总计" (total) sheet ---
$src = $wb.Worksheets.Item("2021-Q4")
$totalBefore = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalBefore)
$q1.Name = "2022-Q1"

# NOTE: inserting a sheet shifts indices, so any handle obtained before the
# Add() that pointed at the "总计" sheet (or anything at/after its position)
# is no longer trustworthy - re-resolve it by name.
$total = $wb.Worksheets.Item("总计")

# Header row - copy number/border/font formatting from an existing quarterly
# sheet's header row, then overwrite the text.
$src.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data row 2 - single fund holding for the new quarter
$src.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)
$q1.Range("A2").Value = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "000522"
$q1.Range("C2").Value = "华润元大信息传媒科技混合"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "1.50"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "70.63"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "3.30"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.0495"
$q1.Range("H2").Value = 9

# --- 2. Insert a new summary row into "总计" for the 2022-Q1 quarter ---
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.05

# Renumber the 0-based index column for the rows that shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

Write-Host "2022-Q1 sheet added and 总计 sheet updated"
